$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.027.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.84%  '
$ws.Range("D3").Value = "'1.547.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").Value = "'286.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").Value = "'0.3811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.19%  '
$ws.Range("D8").Value = "'0.3274"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("D9").Value = "'43.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.87%  '
$ws.Range("D10").Value = "'1.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").Value = "'0.07345"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = "'19.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.06%  '
$ws.Range("D14").Value = "'5.759"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.93%  '
$ws.Range("D15").Value = "'6.712"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.73%  '
$ws.Range("D16").Value = "'1.552.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").Value = "'0.00001072"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.55%  '
$ws.Range("D18").Value = "'0.06634"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'85.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = "'6.329"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = "'15.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").Value = "'11.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").Value = "'22.079.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.61%  '
$ws.Range("D25").Value = "'2.296"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.69%  '
$ws.Range("D26").Value = "'2.489"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.84%  '
$ws.Range("D27").Value = "'150.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").Value = "'19.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("D29").Value = "'4.936"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = "'121.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.12%  '
$ws.Range("D31").Value = "'1.729.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("D32").Value = "'1.070"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("D33").Value = "'5.828"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.71%  '
$ws.Range("D34").Value = "'1.903"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.43%  '
$ws.Range("D35").Value = "'0.08193"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").Value = "'9.245"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.95%  '
$ws.Range("D37").Value = "'0.06257"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.94%  '
$ws.Range("D38").Value = "'0.02308"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.94%  '
$ws.Range("D39").Value = "'5.234"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("D40").Value = "'0.2135"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.17%  '
$ws.Range("D41").Value = "'1.228"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("D42").Value = "'10.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.58%  '
$ws.Range("D43").Value = "'1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = "'0.5999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.23%  '
$ws.Range("D45").Value = "'13.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("D46").Value = "'3.732"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("D47").Value = "'0.5797"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.03%  '
$ws.Range("D48").Value = "'1.964"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.71%  '
$ws.Range("D49").Value = "'121.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.09%  '
$ws.Range("D50").Value = "'1.169"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.37%  '
$ws.Range("D51").Value = "'0.06999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.10%  '